$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the four new "transport" header cells to row 2 (J2:M2). This also
# appends the corresponding four new shared strings, and grows the used
# range / dimension to A1:M2 automatically.
$ws.Range("J2").Value = "运输公司"
$ws.Range("K2").Value = "运输方式"
$ws.Range("L2").Value = "运输单号"
$ws.Range("M2").Value = "收货地址"

# Move the selection to the newly-added last cell, matching the saved
# view state in the workbook.
$ws.Range("M2").Select() | Out-Null
